# Insert a new weekly price-report row for "Femacal de La Calera - Bruselas
# (repollito)" ahead of the existing history (row 64), shifting the rest of
# the rows (old 64..95) down by one (they become 65..96, unchanged), and
# fill the new row with the latest report's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 64:95 down to 65:96, creating a blank row 64.
$ws.Rows("64:64").Insert()

# Populate the new row 64 with the newest report (all the "static" columns
# mirror every other row in this subconjunto; only the date/volume/price
# columns are specific to this report).
$ws.Cells.Item(64, 1).Value = 3
$ws.Cells.Item(64, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(64, 3).Value = "Coquimbo"
$ws.Cells.Item(64, 4).Value = 44839
$ws.Cells.Item(64, 5).Value = 5
$ws.Cells.Item(64, 6).Value = 100112035
$ws.Cells.Item(64, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 40
$ws.Cells.Item(64, 11).Value = 15000
$ws.Cells.Item(64, 12).Value = 15000
$ws.Cells.Item(64, 13).Value = 15000
$ws.Cells.Item(64, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(64, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(64, 16).Value = 1000
$ws.Cells.Item(64, 17).Value = 15
$ws.Cells.Item(64, 18).Value = "Hortaliza"
